# Swap the B and E:AD cell values between two rows for each given pair.
# Columns A (id), C (Div) and D (Date) remain untouched in each row -
# only the match id (B) and the stats/odds columns (E through AD) are
# exchanged between the paired rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, $row1, $row2)

    # Swap column B (match id)
    $b1 = $ws.Cells.Item($row1, 2).Value2
    $b2 = $ws.Cells.Item($row2, 2).Value2
    $ws.Cells.Item($row1, 2).Value2 = $b2
    $ws.Cells.Item($row2, 2).Value2 = $b1

    # Swap columns E (5) through AD (30)
    for ($col = 5; $col -le 30; $col++) {
        $v1 = $ws.Cells.Item($row1, $col).Value2
        $v2 = $ws.Cells.Item($row2, $col).Value2
        $ws.Cells.Item($row1, $col).Value2 = $v2
        $ws.Cells.Item($row2, $col).Value2 = $v1
    }
}

Swap-RowData $ws 110 111
Swap-RowData $ws 237 238
Swap-RowData $ws 249 250
